$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "24÷7=3, 3" "67÷5=13, 2"
Replace-Text "70÷5=14, 0" "18÷3=6, 0"
Replace-Text "47÷5=9, 2" "52÷5=10, 2"
Replace-Text "98÷5=19, 3" "30÷3=10, 0"
Replace-Text "44÷8=5, 4" "87÷7=12, 3"
Replace-Text "96÷7=13, 5" "16÷2=8, 0"
Replace-Text "15÷3=5, 0" "86÷3=28, 2"
Replace-Text "26÷5=5, 1" "87÷3=29, 0"
Replace-Text "41÷7=5, 6" "57÷8=7, 1"
Replace-Text "48÷6=8, 0" "76÷4=19, 0"
Replace-Text "16÷8=2, 0" "17÷5=3, 2"
Replace-Text "80÷3=26, 2" "90÷8=11, 2"
Replace-Text "14÷7=2, 0" "59÷4=14, 3"
Replace-Text "17÷2=8, 1" "75÷3=25, 0"
Replace-Text "28÷2=14, 0" "13÷3=4, 1"
Replace-Text "53÷9=5, 8" "79÷9=8, 7"
Replace-Text "48÷7=6, 6" "96÷3=32, 0"
Replace-Text "54÷7=7, 5" "64÷4=16, 0"
Replace-Text "15÷4=3, 3" "72÷6=12, 0"
Replace-Text "58÷5=11, 3" "42÷9=4, 6"
Replace-Text "76÷8=9, 4" "39÷3=13, 0"
Replace-Text "35÷9=3, 8" "81÷9=9, 0"
Replace-Text "64÷5=12, 4" "24÷4=6, 0"
Replace-Text "73÷8=9, 1" "18÷9=2, 0"
Replace-Text "87÷2=43, 1" "22÷5=4, 2"
